# Auto-generated edit script: updates Leve profit-calculation sheets
# (currentAveragePrice / LevePrice / LeveProfit columns) with refreshed
# market-board data, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 586.8  # H2
$ws.Cells.Item(2, 9).Value = 619.8333  # I2
$ws.Cells.Item(2, 10).Value = 537.25  # J2
$ws.Cells.Item(2, 11).Value = 619.8333  # K2
$ws.Cells.Item(2, 12).Value = 537.25  # L2
$ws.Cells.Item(2, 13).Value = -506.8333  # M2
$ws.Cells.Item(2, 14).Value = -763.25  # N2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 224.66667  # H4
$ws.Cells.Item(4, 9).Value = 224.66667  # I4
$ws.Cells.Item(4, 11).Value = 224.66667  # K4
$ws.Cells.Item(4, 13).Value = -110.66667  # M4

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2998.3333  # H62
$ws.Cells.Item(62, 9).Value = 2998.3333  # I62
$ws.Cells.Item(62, 11).Value = 2998.3333  # K62
$ws.Cells.Item(62, 13).Value = -2374.3333  # M62

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2998.3333  # H65
$ws.Cells.Item(65, 9).Value = 2998.3333  # I65
$ws.Cells.Item(65, 11).Value = 14991.6665  # K65
$ws.Cells.Item(65, 13).Value = -11871.6665  # M65

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 4000  # H70
$ws.Cells.Item(70, 9).Value = 3500  # I70
$ws.Cells.Item(70, 11).Value = 10500  # K70
$ws.Cells.Item(70, 13).Value = -10230  # M70

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 4000  # H73
$ws.Cells.Item(73, 9).Value = 3500  # I73
$ws.Cells.Item(73, 11).Value = 10500  # K73
$ws.Cells.Item(73, 13).Value = -9564  # M73

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 1993  # H113
$ws.Cells.Item(113, 9).Value = 2000  # I113
$ws.Cells.Item(113, 10).Value = 1979  # J113
$ws.Cells.Item(113, 11).Value = 2000  # K113
$ws.Cells.Item(113, 12).Value = 1979  # L113
$ws.Cells.Item(113, 13).Value = 1254  # M113
$ws.Cells.Item(113, 14).Value = -8487  # N113

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 99  # H2
$ws.Cells.Item(2, 9).Value = 99  # I2
$ws.Cells.Item(2, 10).Value = 0  # J2
$ws.Cells.Item(2, 11).Value = 99  # K2
$ws.Cells.Item(2, 12).Value = 0  # L2
$ws.Cells.Item(2, 13).Value = 14  # M2
$ws.Cells.Item(2, 14).ClearContents()  # N2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1010.8182  # H32
$ws.Cells.Item(32, 9).Value = 1010.8182  # I32
$ws.Cells.Item(32, 11).Value = 1010.8182  # K32
$ws.Cells.Item(32, 13).Value = -723.8182  # M32

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 25362.666  # H44
$ws.Cells.Item(44, 9).Value = 3044  # I44
$ws.Cells.Item(44, 10).Value = 70000  # J44
$ws.Cells.Item(44, 11).Value = 3044  # K44
$ws.Cells.Item(44, 12).Value = 70000  # L44
$ws.Cells.Item(44, 13).Value = -2556  # M44
$ws.Cells.Item(44, 14).Value = -70976  # N44

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 4179.143  # H88
$ws.Cells.Item(88, 9).Value = 2851  # I88
$ws.Cells.Item(88, 11).Value = 2851  # K88
$ws.Cells.Item(88, 13).Value = -2445  # M88

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 4179.143  # H91
$ws.Cells.Item(91, 9).Value = 2851  # I91
$ws.Cells.Item(91, 11).Value = 2851  # K91
$ws.Cells.Item(91, 13).Value = -1447  # M91

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(95, 8).Value = 28613.4  # H95
$ws.Cells.Item(95, 10).Value = 28613.4  # J95
$ws.Cells.Item(95, 12).Value = 28613.4  # L95
$ws.Cells.Item(95, 14).Value = -34105.4  # N95

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 41668970  # H97
$ws.Cells.Item(97, 9).Value = 66669750  # I97
$ws.Cells.Item(97, 11).Value = 66669750  # K97
$ws.Cells.Item(97, 13).Value = -66669254  # M97

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 11107790  # H102
$ws.Cells.Item(102, 9).Value = 848309.0600000001  # I102
$ws.Cells.Item(102, 11).Value = 848309.0600000001  # K102
$ws.Cells.Item(102, 13).Value = -846687.0600000001  # M102

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(103, 8).Value = 29069  # H103
$ws.Cells.Item(103, 10).Value = 29069  # J103
$ws.Cells.Item(103, 12).Value = 29069  # L103
$ws.Cells.Item(103, 14).Value = -31413  # N103

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 99  # H116
$ws.Cells.Item(116, 9).Value = 99  # I116
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 11).Value = 99  # K116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 13).Value = 2195  # M116
$ws.Cells.Item(116, 14).ClearContents()  # N116

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 99  # H3
$ws.Cells.Item(3, 9).Value = 99  # I3
$ws.Cells.Item(3, 10).Value = 0  # J3
$ws.Cells.Item(3, 11).Value = 99  # K3
$ws.Cells.Item(3, 12).Value = 0  # L3
$ws.Cells.Item(3, 13).Value = 15  # M3
$ws.Cells.Item(3, 14).ClearContents()  # N3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 93082  # H94
$ws.Cells.Item(94, 9).Value = 123559.336  # I94
$ws.Cells.Item(94, 10).Value = 1650  # J94
$ws.Cells.Item(94, 11).Value = 123559.336  # K94
$ws.Cells.Item(94, 12).Value = 1650  # L94
$ws.Cells.Item(94, 13).Value = -123108.336  # M94
$ws.Cells.Item(94, 14).Value = -2552  # N94

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value = 23727.8  # H106
$ws.Cells.Item(106, 10).Value = 23727.8  # J106
$ws.Cells.Item(106, 12).Value = 23727.8  # L106
$ws.Cells.Item(106, 14).Value = -26251.8  # N106

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 37317.184  # H107
$ws.Cells.Item(107, 9).Value = 37317.184  # I107
$ws.Cells.Item(107, 11).Value = 37317.184  # K107
$ws.Cells.Item(107, 13).Value = -35397.184  # M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 189.54546  # H7
$ws.Cells.Item(7, 9).Value = 213.21053  # I7
$ws.Cells.Item(7, 10).Value = 39.666668  # J7
$ws.Cells.Item(7, 11).Value = 213.21053  # K7
$ws.Cells.Item(7, 12).Value = 39.666668  # L7
$ws.Cells.Item(7, 13).Value = -100.21053  # M7
$ws.Cells.Item(7, 14).Value = -265.666668  # N7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 18750  # H28
$ws.Cells.Item(28, 10).Value = 18750  # J28
$ws.Cells.Item(28, 12).Value = 18750  # L28
$ws.Cells.Item(28, 14).Value = -19240  # N28

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 2639.6  # H35
$ws.Cells.Item(35, 9).Value = 2639.6  # I35
$ws.Cells.Item(35, 10).Value = 0  # J35
$ws.Cells.Item(35, 11).Value = 2639.6  # K35
$ws.Cells.Item(35, 12).Value = 0  # L35
$ws.Cells.Item(35, 13).Value = -2345.6  # M35
$ws.Cells.Item(35, 14).ClearContents()  # N35

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1798.6666  # H94
$ws.Cells.Item(94, 9).Value = 1697  # I94
$ws.Cells.Item(94, 10).Value = 1849.5  # J94
$ws.Cells.Item(94, 11).Value = 1697  # K94
$ws.Cells.Item(94, 12).Value = 1849.5  # L94
$ws.Cells.Item(94, 13).Value = -1246  # M94
$ws.Cells.Item(94, 14).Value = -2751.5  # N94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 560666  # H141
$ws.Cells.Item(141, 9).Value = 14500  # I141
$ws.Cells.Item(141, 11).Value = 14500  # K141
$ws.Cells.Item(141, 13).Value = -9320  # M141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 369.3846  # H2
$ws.Cells.Item(2, 9).Value = 396  # I2
$ws.Cells.Item(2, 11).Value = 2376  # K2
$ws.Cells.Item(2, 13).Value = -2263  # M2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 963.8461  # H32
$ws.Cells.Item(32, 9).Value = 91.25  # I32
$ws.Cells.Item(32, 10).Value = 2360  # J32
$ws.Cells.Item(32, 11).Value = 273.75  # K32
$ws.Cells.Item(32, 12).Value = 7080  # L32
$ws.Cells.Item(32, 13).Value = 9.25  # M32
$ws.Cells.Item(32, 14).Value = -7646  # N32

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 168  # H41
$ws.Cells.Item(41, 9).Value = 168  # I41
$ws.Cells.Item(41, 10).Value = 0  # J41
$ws.Cells.Item(41, 11).Value = 504  # K41
$ws.Cells.Item(41, 12).Value = 0  # L41
$ws.Cells.Item(41, 13).Value = -166  # M41
$ws.Cells.Item(41, 14).ClearContents()  # N41

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 0  # H109
$ws.Cells.Item(109, 9).Value = 0  # I109
$ws.Cells.Item(109, 10).Value = 0  # J109
$ws.Cells.Item(109, 11).Value = 0  # K109
$ws.Cells.Item(109, 12).Value = 0  # L109
$ws.Cells.Item(109, 14).ClearContents()  # N109

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 2000  # H140
$ws.Cells.Item(140, 9).Value = 2000  # I140
$ws.Cells.Item(140, 11).Value = 6000  # K140
$ws.Cells.Item(140, 13).Value = -820  # M140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2442.2856  # H80
$ws.Cells.Item(80, 9).Value = 3130  # I80
$ws.Cells.Item(80, 10).Value = 1926.5  # J80
$ws.Cells.Item(80, 11).Value = 3130  # K80
$ws.Cells.Item(80, 12).Value = 1926.5  # L80
$ws.Cells.Item(80, 13).Value = -2132  # M80
$ws.Cells.Item(80, 14).Value = -3922.5  # N80

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2442.2856  # H83
$ws.Cells.Item(83, 9).Value = 3130  # I83
$ws.Cells.Item(83, 10).Value = 1926.5  # J83
$ws.Cells.Item(83, 11).Value = 15650  # K83
$ws.Cells.Item(83, 12).Value = 9632.5  # L83
$ws.Cells.Item(83, 13).Value = -10658  # M83
$ws.Cells.Item(83, 14).Value = -19616.5  # N83

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 10608  # H98
$ws.Cells.Item(98, 10).Value = 10608  # J98
$ws.Cells.Item(98, 12).Value = 10608  # L98
$ws.Cells.Item(98, 14).Value = -16598  # N98

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(105, 8).Value = 12300  # H105
$ws.Cells.Item(105, 10).Value = 12300  # J105
$ws.Cells.Item(105, 12).Value = 12300  # L105
$ws.Cells.Item(105, 14).Value = -19288  # N105

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 0  # H113
$ws.Cells.Item(113, 9).Value = 0  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 0  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).ClearContents()  # M113

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7930  # H122
$ws.Cells.Item(122, 10).Value = 8999.75  # J122
$ws.Cells.Item(122, 12).Value = 26999.25  # L122
$ws.Cells.Item(122, 14).Value = -31899.25  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3497.375  # H22
$ws.Cells.Item(22, 9).Value = 3245.8  # I22
$ws.Cells.Item(22, 11).Value = 3245.8  # K22
$ws.Cells.Item(22, 13).Value = -2950.8  # M22

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3497.375  # H27
$ws.Cells.Item(27, 9).Value = 3245.8  # I27
$ws.Cells.Item(27, 11).Value = 3245.8  # K27
$ws.Cells.Item(27, 13).Value = -3138.8  # M27

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(39, 8).Value = 27497.5  # H39
$ws.Cells.Item(39, 10).Value = 35000  # J39
$ws.Cells.Item(39, 12).Value = 35000  # L39
$ws.Cells.Item(39, 14).Value = -35920  # N39

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(54, 8).Value = 42509.5  # H54
$ws.Cells.Item(54, 10).Value = 42509.5  # J54
$ws.Cells.Item(54, 12).Value = 42509.5  # L54
$ws.Cells.Item(54, 14).Value = -43797.5  # N54

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 8642.5  # H68
$ws.Cells.Item(68, 9).Value = 8642.5  # I68
$ws.Cells.Item(68, 11).Value = 8642.5  # K68
$ws.Cells.Item(68, 13).Value = -7893.5  # M68

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 8642.5  # H71
$ws.Cells.Item(71, 9).Value = 8642.5  # I71
$ws.Cells.Item(71, 11).Value = 43212.5  # K71
$ws.Cells.Item(71, 13).Value = -39468.5  # M71

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 66667610  # H93
$ws.Cells.Item(93, 9).Value = 83334260  # I93
$ws.Cells.Item(93, 11).Value = 83334260  # K93
$ws.Cells.Item(93, 13).Value = -83333012  # M93

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(104, 8).Value = 27045.8  # H104
$ws.Cells.Item(104, 10).Value = 27633.111  # J104
$ws.Cells.Item(104, 12).Value = 27633.111  # L104
$ws.Cells.Item(104, 14).Value = -34621.111  # N104

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 10229.25  # H106
$ws.Cells.Item(106, 10).Value = 10229.25  # J106
$ws.Cells.Item(106, 12).Value = 10229.25  # L106
$ws.Cells.Item(106, 14).Value = -12753.25  # N106

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3493.2415  # H122
$ws.Cells.Item(122, 9).Value = 2800.2856  # I122
$ws.Cells.Item(122, 11).Value = 8400.856800000001  # K122
$ws.Cells.Item(122, 13).Value = -5950.856800000001  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 0  # H54
$ws.Cells.Item(54, 9).Value = 0  # I54
$ws.Cells.Item(54, 10).Value = 0  # J54
$ws.Cells.Item(54, 11).Value = 0  # K54
$ws.Cells.Item(54, 12).Value = 0  # L54
$ws.Cells.Item(54, 14).ClearContents()  # N54
